$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2, 3, 0.04969111061828357),
    @(2, 4, 0.2107476850596086),
    @(2, 5, 0.1680132064088014),
    @(2, 6, 1.276671991392604),
    @(2, 7, 0.7614461594882869),
    @(2, 8, 0.7362247279386338),
    @(2, 10, 0.1838998710621524),
    @(2, 14, 2.873581687696344),
    @(2, 15, 3.002323868125757),
    @(3, 3, 0.0440893987137656),
    @(3, 4, 0.2101148881431385),
    @(3, 5, 0.1653876853705007),
    @(3, 6, 1.235354089528187),
    @(3, 7, 0.7201967886403509),
    @(3, 8, 0.7216148393470121),
    @(3, 10, 0.1786267740460445),
    @(3, 14, 2.562605684679454),
    @(3, 15, 2.883864273825111),
    @(4, 3, 0.04066648080491575),
    @(4, 4, 0.2098108844849804),
    @(4, 5, 0.1638661144860727),
    @(4, 6, 1.210749221148646),
    @(4, 7, 0.6953400628314057),
    @(4, 8, 0.7130539717894919),
    @(4, 10, 0.1755003753982649),
    @(4, 14, 2.371325805375761),
    @(4, 15, 2.812967843030833),
    @(5, 3, 0.03927570747224252),
    @(5, 4, 0.2097082792863816),
    @(5, 5, 0.1632687989182777),
    @(5, 6, 1.200914363468527),
    @(5, 7, 0.6853285548016999),
    @(5, 8, 0.7096682597508703),
    @(5, 10, 0.1742542540487335),
    @(5, 14, 2.293303068607429),
    @(5, 15, 2.784537893400397),
    @(6, 3, 0.03904501583026843),
    @(6, 4, 0.2096925275952657),
    @(6, 5, 0.1631709877360805),
    @(6, 6, 1.199292866790799),
    @(6, 7, 0.6836732519807072),
    @(6, 8, 0.7091122775099876),
    @(6, 10, 0.1740490201064162),
    @(6, 14, 2.280343261403573),
    @(6, 15, 2.779844922308712),
    @(7, 3, 0.04064770786656879),
    @(7, 4, 0.2098094145296088),
    @(7, 5, 0.1638579668473099),
    @(7, 6, 1.21061580867098),
    @(7, 7, 0.6952045678091281),
    @(7, 8, 0.7130078943511364),
    @(7, 10, 0.1754834568576271),
    @(7, 14, 2.370273851395496),
    @(7, 15, 2.812582562290402),
    @(8, 3, 0.04775616062549659),
    @(8, 4, 0.2105119633123351),
    @(8, 5, 0.1670891265871504),
    @(8, 6, 1.262266643986948),
    @(8, 7, 0.7471253876036883),
    @(8, 8, 0.7311021441246055),
    @(8, 10, 0.1820585575030407),
    @(8, 14, 2.766433886209654),
    @(8, 15, 2.96109663284426),
    @(9, 3, 0.0618310090167995),
    @(9, 4, 0.2125597742668219),
    @(9, 5, 0.1741451481489946),
    @(9, 6, 1.369645234141316),
    @(9, 7, 0.8527049597044538),
    @(9, 8, 0.7698429204376112),
    @(9, 10, 0.1958398391215894),
    @(9, 14, 3.540180268007646),
    @(9, 15, 3.266996416250606),
    @(10, 3, 0.0722605177838318),
    @(10, 4, 0.2144724038754475),
    @(10, 5, 0.1797709608777254),
    @(10, 6, 1.452295328595596),
    @(10, 7, 0.9326177770774677),
    @(10, 8, 0.8003065429869878),
    @(10, 10, 0.2065133899798752),
    @(10, 14, 4.10623028343673),
    @(10, 15, 3.500816883837047),
    @(11, 3, 0.07702584964360426),
    @(11, 4, 0.2154310724398982),
    @(11, 5, 0.1824269134106373),
    @(11, 6, 1.490721737755166),
    @(11, 7, 0.9694922927976108),
    @(11, 8, 0.8146031334345594),
    @(11, 10, 0.2114898347396519),
    @(11, 14, 4.363110593465422),
    @(11, 15, 3.60918982072053),
    @(12, 3, 0.07883345327934421),
    @(12, 4, 0.2158068242922155),
    @(12, 5, 0.183446606245603),
    @(12, 6, 1.505392597165411),
    @(12, 7, 0.9835314671214803),
    @(12, 8, 0.8200801398335216),
    @(12, 10, 0.2133917932220157),
    @(12, 14, 4.460285735714251),
    @(12, 15, 3.650518376121966),
    @(13, 3, 0.07844401556036473),
    @(13, 4, 0.2157253337201013),
    @(13, 5, 0.1832263766265498),
    @(13, 6, 1.502227637601607),
    @(13, 7, 0.9805045085584823),
    @(13, 8, 0.818897755549358),
    @(13, 10, 0.2129813932859292),
    @(13, 14, 4.439361943450422),
    @(13, 15, 3.641604607028398),
    @(14, 3, 0.0771745002290487),
    @(14, 4, 0.2154617308372053),
    @(14, 5, 0.1825105245772747),
    @(14, 6, 1.491926318211313),
    @(14, 7, 0.9706457858187889),
    @(14, 8, 0.8150524626724689),
    @(14, 10, 0.2116459590375399),
    @(14, 14, 4.371107314139522),
    @(14, 15, 3.612584123659303),
    @(15, 3, 0.07639728811486179),
    @(15, 4, 0.2153019231388527),
    @(15, 5, 0.1820738609860157),
    @(15, 6, 1.485632054026979),
    @(15, 7, 0.9646168950917229),
    @(15, 8, 0.8127053441018575),
    @(15, 10, 0.2108302469391248),
    @(15, 14, 4.329286057409945),
    @(15, 15, 3.594846060814973),
    @(16, 3, 0.07194951869236377),
    @(16, 4, 0.2144115341665582),
    @(16, 5, 0.1795993373709948),
    @(16, 6, 1.449800784119205),
    @(16, 7, 0.9302184883761697),
    @(16, 8, 0.7993810632741543),
    @(16, 10, 0.2061906098194157),
    @(16, 14, 4.089429168003562),
    @(16, 15, 3.493774967663057),
    @(17, 3, 0.06922635870277816),
    @(17, 4, 0.2138879912362484),
    @(17, 5, 0.1781060963176273),
    @(17, 6, 1.428031990147105),
    @(17, 7, 0.9092502194449139),
    @(17, 8, 0.7913194645214787),
    @(17, 10, 0.2033753944021726),
    @(17, 14, 3.94211849063862),
    @(17, 15, 3.4322861830546),
    @(18, 3, 0.06766203167551055),
    @(18, 4, 0.2135952039426456),
    @(18, 5, 0.177256327565317),
    @(18, 6, 1.415589105775723),
    @(18, 7, 0.8972389035548076),
    @(18, 8, 0.7867239148270926),
    @(18, 10, 0.2017675400542345),
    @(18, 14, 3.857331695637754),
    @(18, 15, 3.397108262714426),
    @(19, 3, 0.06713271233556384),
    @(19, 4, 0.213497504228684),
    @(19, 5, 0.176970173187982),
    @(19, 6, 1.411389532417957),
    @(19, 7, 0.8931804903969294),
    @(19, 8, 0.7851750244728919),
    @(19, 10, 0.2012251005842671),
    @(19, 14, 3.828614786364199),
    @(19, 15, 3.385230004966104),
    @(20, 3, 0.06951604009691437),
    @(20, 4, 0.2139428601582267),
    @(20, 5, 0.1782641119323145),
    @(20, 6, 1.430341244444293),
    @(20, 7, 0.9114772440300953),
    @(20, 8, 0.7921733633517647),
    @(20, 10, 0.2036739001438832),
    @(20, 14, 3.957806003280837),
    @(20, 15, 3.438812218089822),
    @(21, 3, 0.07754730360997542),
    @(21, 4, 0.2155388121715305),
    @(21, 5, 0.1827204090742356),
    @(21, 6, 1.494948815327632),
    @(21, 7, 0.9735394754994786),
    @(21, 8, 0.8161802021215863),
    @(21, 10, 0.2120377331518881),
    @(21, 14, 4.391158149571254),
    @(21, 15, 3.621100258499951),
    @(22, 3, 0.08281416027682553),
    @(22, 4, 0.21665602212056),
    @(22, 5, 0.1857141314301884),
    @(22, 6, 1.537871016488282),
    @(22, 7, 1.014541625562316),
    @(22, 8, 0.8322384904946318),
    @(22, 10, 0.2176059498936098),
    @(22, 14, 4.673791817957863),
    @(22, 15, 3.741927523114384),
    @(23, 3, 0.08000147273300229),
    @(23, 4, 0.2160529654793635),
    @(23, 5, 0.1841088788204104),
    @(23, 6, 1.514898655921542),
    @(23, 7, 0.9926174766382303),
    @(23, 8, 0.8236341221845009),
    @(23, 10, 0.2146247283206009),
    @(23, 14, 4.523002190005627),
    @(23, 15, 3.677284467669892),
    @(24, 3, 0.06938507123904003),
    @(24, 4, 0.2139180283619169),
    @(24, 5, 0.1781926459288243),
    @(24, 6, 1.429297005451531),
    @(24, 7, 0.910470270384053),
    @(24, 8, 0.7917871936849394),
    @(24, 10, 0.2035389125315561),
    @(24, 14, 3.950713976768498),
    @(24, 15, 3.43586125950452),
    @(25, 3, 0.05800826539810089),
    @(25, 4, 0.2119340913195487),
    @(25, 5, 0.1721589486874109),
    @(25, 6, 1.339939663534139),
    @(25, 7, 0.82373493451135),
    @(25, 8, 0.759012223452487),
    @(25, 10, 0.1920158612543617),
    @(25, 14, 3.331249627311365),
    @(25, 15, 3.182658497228545)
)

foreach ($row in $data) {
    $r = $row[0]
    $c = $row[1]
    $v = $row[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Host "Applied $($data.Count) cell updates"
